$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: new data row appended below existing header (row1) and data (row2)

# A3 holds a date-like text value "09/13/2025". Assigning it directly would make
# Excel auto-convert it into a date serial number, so we briefly force a text
# number format, assign the value, then clear the format override so the cell
# ends up with the default style (matching the plain, unstyled text cells used
# elsewhere in the sheet) while keeping the value stored as text.
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "09/13/2025"
$ws.Range("A3").ClearFormats()

$ws.Range("B3").Value = "DeepSeekV3"

$ws.Range("C3").Value = 0.4957627118644068
$ws.Range("D3").Value = 0.4034482758620689
$ws.Range("E3").Value = 0.4448669201520912
$ws.Range("F3").Value = 0.5688541426844711
$ws.Range("G3").Value = 0.3918825071348214
$ws.Range("H3").Value = 0.4369894502903656
$ws.Range("I3").Value = 0.58283460620904
$ws.Range("J3").Value = 0.4034482758620689
$ws.Range("K3").Value = 0.4565638851009324
$ws.Range("L3").Value = 0.5635593220338984
$ws.Range("M3").Value = 0.4602076124567474
$ws.Range("N3").Value = 0.5066666666666667
$ws.Range("O3").Value = 133
$ws.Range("P3").Value = 103
$ws.Range("Q3").Value = 156
$ws.Range("R3").Value = 290
$ws.Range("S3").Value = 0.9558841500672125

$ws.Range("T3").Value = "/home/s27mhusa_hpc/Master-Thesis/Evaluation_Results/Final_TestFiles_12thSeptember_FewShotTest_Broad_Count/ner_evaluation_results_DeepSeekV3_4_shot.txt"
$ws.Range("U3").Value = "/home/s27mhusa_hpc/Master-Thesis/Evaluation_Results/Final_TestFiles_12thSeptember_FewShotTest_Broad_Count/Stats/ner_evaluation_stats_DeepSeekV3_4_shot.txt"
$ws.Range("V3").Value = "4 MLGPU"
$ws.Range("W3").Value = "0.003 kWh"
